$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Calculations sheet: insert a new blank row above the old row 10 (the "lb
# per metric ton" conversion block), shifting everything below down by one.
# This turns the old A10/A11/B11 ("lb per metric ton" / "2204.62" /
# "lb/metric ton") into A11/A12/B12, the old A13 label into A14, and the
# old A14 formula into A15 - and Excel auto-updates the formula's relative
# reference (A11 -> A12) plus every downstream reference to
# Calculations!$A$14 (on the EHPpUC sheet) to Calculations!$A$15.
# ---------------------------------------------------------------------------
$wsCalc = $wb.Worksheets.Item("Calculations")
$wsCalc.Rows.Item(9).Insert()

# Updated hydrogen energy density source value (BTU/lb), re-styled to
# Arial 12pt black to match the other cited-source values in the sheet.
$wsCalc.Range("A8").Value = 60920
$wsCalc.Range("A8").Font.Name = "Arial"
$wsCalc.Range("A8").Font.Size = 12
$wsCalc.Range("A8").Font.Color = 0
$wsCalc.Rows.Item(8).RowHeight = 15.4

# New reviewer notes placed beside the "lb per metric ton" conversion block.
$wsCalc.Range("C12").Value = "The study they are citing uses 3 scenarios of Fuel Cell Electric Vehicle adoption. Then it calculates the amount of hydrogen needed to suppor those vehicles. Then it calculates the electrolyzer capacity needed to supply that hydrogen. "
$wsCalc.Range("C13").Value = "So, I think it's fair, using the EPS assumptions of 24/7/365 operation, that smallest electrolyzer you would need to produce 1.39e10 annual Btu would be 1 MW."
$wsCalc.Range("C14").Value = "No reason to think this would be different for Texas."

# Give the new note cells the same look as the other highlighted note cells
# used elsewhere in the workbook (Calibri, accent5 theme color).
$wsCalc.Range("C12:C14").Font.ThemeColor = 9

# Restore the selected cell on each sheet to match the saved workbook state.
$wsCalc.Range("J6").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B6").Select()

$wsEHPpUC = $wb.Worksheets.Item("EHPpUC")
$wsEHPpUC.Range("B2").Select()

$wsAbout.Activate()
